# Remove the "Proposed incoming transfer name" / "[SponsorName]" row
# from the summary table (the second row of the second table in the
# document, right after the "Project name" / "[ProjectName]" row).

$d = $word.ActiveDocument

$targetTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if ($t.Range.Text -like "*Proposed incoming transfer*SponsorName*") {
        $targetTable = $t
        break
    }
}

for ($r = $targetTable.Rows.Count; $r -ge 1; $r--) {
    $row = $targetTable.Rows.Item($r)
    if ($row.Range.Text -like "*Proposed incoming transfer*") {
        $row.Delete()
    }
}
